# UPDATE BUDGET TO FIT RANGE
# Extends the Budget sheet's data range from A1:B19 to A1:B49, applying a
# 5% uplift to the existing forecast rows (14-19) and carrying the same
# quarterly-step growth pattern forward through FY25 (rows 20-49). The
# final 12 rows (38-49) are displayed with an extra decimal of precision
# via a new currency number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Bump the already-forecast rows (14:19) by 5%
# ---------------------------------------------------------------------
$ws.Range("B14").Value = 840000
$ws.Range("B15").Value = 840000
$ws.Range("B16").Value = 1050000
$ws.Range("B17").Value = 1050000
$ws.Range("B18").Value = 1155000
$ws.Range("B19").Value = 1155000

# ---------------------------------------------------------------------
# 2. Stamp the date/currency formatting of row 19 across the new rows
#    (20:49) first, so the new cells inherit the existing "Date" (s=3)
#    and "Currency" (s=2) styles instead of landing on General.
# ---------------------------------------------------------------------
$ws.Range("A19:B19").Copy()
$ws.Range("A20:B49").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Write the new date / budget values for rows 20:49
# ---------------------------------------------------------------------
$newRows = @(
    @(20, 44378, 1575000),
    @(21, 44409, 1575000),
    @(22, 44440, 1575000),
    @(23, 44470, 1575000),
    @(24, 44501, 1575000),
    @(25, 44531, 2100000),
    @(26, 44562, 882000),
    @(27, 44593, 882000),
    @(28, 44621, 1102500),
    @(29, 44652, 1102500),
    @(30, 44682, 1212750),
    @(31, 44713, 1212750),
    @(32, 44743, 1653750),
    @(33, 44774, 1653750),
    @(34, 44805, 1653750),
    @(35, 44835, 1653750),
    @(36, 44866, 1653750),
    @(37, 44896, 2205000),
    @(38, 44927, 904049.99999999988),
    @(39, 44958, 904049.99999999988),
    @(40, 44986, 1130062.5),
    @(41, 45017, 1130062.5),
    @(42, 45047, 1243068.75),
    @(43, 45078, 1243068.75),
    @(44, 45108, 1695093.7499999998),
    @(45, 45139, 1695093.7499999998),
    @(46, 45170, 1695093.7499999998),
    @(47, 45200, 1695093.7499999998),
    @(48, 45231, 1695093.7499999998),
    @(49, 45261, 2260125)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $dateSerial = $row[1]
    $budget = $row[2]
    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r, 2).Value = $budget
}

# Row 38's budget is actually a formula off row 26 (the first of the new
# 5%-uplifted quarters); re-enter it as a live formula so the cached
# value is derived rather than just a literal.
$ws.Range("B38").Formula = "=1.025*B26"

# ---------------------------------------------------------------------
# 4. Rows 38:49 get a new number format with three decimal places
#    instead of the usual two (new custom numFmt + cellXf).
# ---------------------------------------------------------------------
$ws.Range("B38:B49").NumberFormat = '_("$"* #,##0.000_);_("$"* \(#,##0.000\);_("$"* "-"???_);_(@_)'

# ---------------------------------------------------------------------
# 5. Misc view state to mirror the authoring session
# ---------------------------------------------------------------------
$ws.Range("H8").Select() | Out-Null
